$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$t98 = @'
Background
id="Par1">The global numbers of confirmed cases and deceased critically ill patients with COVID-19 are increasing.

 However, the clinical course, and the 60-day mortality and its predictors in critically ill patients have not been fully elucidated.

 The aim of this study is to identify the clinical course, and 60-day mortality and its predictors in critically ill patients with COVID-19.
Methods
id="Par2">Critically ill adult patients admitted to intensive care units (ICUs) from 3 hospitals in Wuhan, China, were included.

 Data on demographic information, preexisting comorbidities, laboratory findings at ICU admission, treatments, clinical outcomes, and results of SARS-CoV-2 RNA tests and of serum SARS-CoV-2 IgM were collected including the duration between symptom onset and negative conversion of SARS-CoV-2 RNA.


Results
id="Par3">Of 1748 patients with COVID-19, 239 (13.7%) critically ill patients were included.

 Complications included acute respiratory distress syndrome (ARDS) in 164 (68.6%) patients, coagulopathy in 150 (62.7%) patients, acute cardiac injury in 103 (43.1%) patients, and acute kidney injury (AKI) in 119 (49.8%) patients, which occurred 15.5 days, 17 days, 18.5 days, and 19 days after the symptom onset, respectively.

 The median duration of the negative conversion of SARS-CoV-2 RNA was 30 (range 6–81) days in 49 critically ill survivors that were identified.

 A total of 147 (61.5%) patients deceased by 60 days after ICU admission.

 The median duration between ICU admission and decease was 12 (range 3–36).

 Cox proportional-hazards regression analysis revealed that age older than 65 years, thrombocytopenia at ICU admission, ARDS, and AKI independently predicted the 60-day mortality.


Conclusions
id="Par4">Severe complications are common and the 60-day mortality of critically ill patients with COVID-19 is considerably high.

 The duration of the negative conversion of SARS-CoV-2 RNA and its association with the severity of critically ill patients with COVID-19 should be seriously considered and further studied.



'@
$t99 = @'
[Jiqian%Xu%NULL%1, Xiaobo%Yang%NULL%0, Luyu%Yang%NULL%1, Xiaojing%Zou%NULL%1, Yaxin%Wang%NULL%1, Yongran%Wu%NULL%1, Ting%Zhou%NULL%1, Yin%Yuan%NULL%1, Hong%Qi%NULL%1, Shouzhi%Fu%NULL%1, Hong%Liu%NULL%1, Jia’an%Xia%NULL%1, Zhengqin%Xu%NULL%1, Yuan%Yu%NULL%1, Ruiting%Li%NULL%1, Yaqi%Ouyang%NULL%1, Rui%Wang%NULL%1, Lehao%Ren%NULL%1, Yingying%Hu%NULL%1, Dan%Xu%NULL%1, Xin%Zhao%NULL%1, Shiying%Yuan%shiying_yuan@163.com%1, Dingyu%Zhang%1813886398@qq.com%0, You%Shang%you_shanghust@163.com%1]
'@
$t136 = @'
Background
id="Par1">It had been shown that High-flow nasal cannula (HFNC) is an effective initial support strategy for patients with acute respiratory failure.

 However, the efficacy of HFNC for patients with COVID-19 has not been established.

 This study was performed to assess the efficacy of HFNC for patients with COVID-19 and describe early predictors of HFNC treatment success in order to develop a prediction tool that accurately identifies the need for upgrade respiratory support therapy.


Methods
id="Par2">We retrospectively reviewed the medical records of patients with COVID-19 treated by HFNC in respiratory wards of 2 hospitals in Wuhan between 1 January and 1 March 2020. Overall clinical outcomes, the success rate of HFNC strategy and related respiratory variables were evaluated.


Results
id="Par3">A total of 105 patients were analyzed.

 Of these, 65 patients (61.9%) showed improved oxygenation and were successfully withdrawn from HFNC.

 The PaO2/FiO2 ratio, SpO2/FiO2 ratio and ROX index (SpO2/FiO2*RR) at 6h, 12h and 24h of HFNC initiation were closely related to the prognosis.

 The ROX index after 6h of HFNC initiation (AUROC, 0.798) had good predictive capacity for outcomes of HFNC.

 In the multivariate logistic regression analysis, young age, gender of female, and lower SOFA score all have predictive value, while a ROX index greater than 5.55 at 6 h after initiation was significantly associated with HFNC success (OR, 17.821; 95% CI, 3.741-84.903 p&lt;0.001).


Conclusions
id="Par4">Our study indicated that HFNC was an effective way of respiratory support in the treatment of COVID-19 patients.

 The ROX index after 6h after initiating HFNC had good predictive capacity for HFNC outcomes.



'@
$t137 = @'
[
Ming%Hu%NULL%1, Qiang%Zhou%NULL%1, Ruiqiang%Zheng%NULL%1, Xuyan%Li%NULL%1, Jianmin%Ling%NULL%1, Yumei%Chen%NULL%1, Jing%Jia%NULL%1, Cuihong%Xie%xiecuihong08@163.com%1]
'@

$ws.Range("D16").Value = $t98
$ws.Range("E16").Value = $t99
$ws.Range("D25").Value = $t136
$ws.Range("E25").Value = $t137

